$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.346.09'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '3.712.34'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '403.30'
$ws.Range('E5').Value = '  -5.08%  '
$ws.Range('D6').Value = '127.18'
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('D7').Value = '3.700.65'
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  -5.82%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '0.711'
$ws.Range('E10').Value = '  -6.64%  '
$ws.Range('D11').Value = '0.163'
$ws.Range('E11').Value = '  -10.55%  '
$ws.Range('D12').Value = '0.0000346'
$ws.Range('E12').Value = '  -10.08%  '
$ws.Range('D13').Value = '40.15'
$ws.Range('E13').Value = '  -5.50%  '
$ws.Range('D14').Value = '4.283.89'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '9.55'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '14.40'
$ws.Range('E16').Value = '  +9.41%  '
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').Value = '3.716.02'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '19.32'
$ws.Range('E19').Value = '  -6.38%  '
$ws.Range('D20').Value = '65.522.15'
$ws.Range('E20').Value = '  -1.76%  '
$ws.Range('D21').Value = '1.05'
$ws.Range('E21').Value = '  -6.48%  '
$ws.Range('D22').Value = '404.44'
$ws.Range('E22').Value = '  -9.58%  '
$ws.Range('D23').Value = '14.32'
$ws.Range('E23').Value = '  -7.77%  '
$ws.Range('D24').Value = '84.52'
$ws.Range('E24').Value = '  -4.17%  '
$ws.Range('D25').Value = '3.00'
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('D26').Value = '35.96'
$ws.Range('E26').Value = '  -5.09%  '
$ws.Range('D27').Value = '5.48'
$ws.Range('E27').Value = '  +9.90%  '
$ws.Range('D28').Value = '3.06'
$ws.Range('E28').Value = '  -9.15%  '
$ws.Range('D29').Value = '9.08'
$ws.Range('E29').Value = '  -11.50%  '
$ws.Range('D30').Value = '12.34'
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').Value = '7.03'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('D34').Value = '0.154'
$ws.Range('E34').Value = '  -4.71%  '
$ws.Range('D35').Value = '38.11'
$ws.Range('E35').Value = '  -9.12%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').Value = '55.07'
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0720'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0449'
$ws.Range('E39').Value = '  -8.02%  '
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('D41').Value = '2.77'
$ws.Range('E41').Value = '  -7.76%  '
$ws.Range('D42').Value = '0.133'
$ws.Range('E42').Value = '  -8.38%  '
$ws.Range('D43').Value = '3.14'
$ws.Range('E43').Value = '  +19.95%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '144.56'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '26.28'
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = '3.19'
$ws.Range('E46').Value = '  -8.06%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '2.02'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '2.78'
$ws.Range('E48').Value = '  -5.00%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '4.21'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = '2.53'
$ws.Range('E50').Value = '  -3.80%  '
$ws.Range('D51').Value = '0.287'
$ws.Range('E51').Value = '  -6.12%  '
